$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A28").Value = "Finansielle foretak"
$ws.Range("B28").Value = "Alternative investeringsfond (AIF) utenom verdipapirfond"

$wb.Save()
